$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45171 -> 2023-09-02).
# The update bumps every data row (C2:C97) forward by one day to 45172 (2023-09-03).
for ($row = 2; $row -le 97; $row++) {
    $ws.Cells.Item($row, 3).Value = 45172
}
